$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.943.30"
$ws.Range("E2").Value = "  -1.48%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.280.73"
$ws.Range("E3").Value = "  -2.83%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.51"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.46"
$ws.Range("E6").Value = "  -4.92%  "
$ws.Range("E7").Value = "  -1.04%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.601"
$ws.Range("E9").Value = "  -2.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.90"
$ws.Range("E10").Value = "  -5.45%  "
$ws.Range("E11").Value = "  -2.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.28"
$ws.Range("E12").Value = "  -2.78%  "
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.960"
$ws.Range("E14").Value = "  -3.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.23"
$ws.Range("E15").Value = "  -4.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.626.63"
$ws.Range("E16").Value = "  -2.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.286.59"
$ws.Range("E17").Value = "  -2.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.797.74"
$ws.Range("E18").Value = "  -1.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.54"
$ws.Range("E19").Value = "  -0.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000105"
$ws.Range("E20").Value = "  -0.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "285.42"
$ws.Range("E21").Value = "  +11.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.62"
$ws.Range("E22").Value = "  -3.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.55"
$ws.Range("E23").Value = "  -0.69%  "
$ws.Range("E24").Value = "  -1.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.92"
$ws.Range("E25").Value = "  +6.16%  "
$ws.Range("E26").Value = "  +0.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.75"
$ws.Range("E27").Value = "  -5.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.28"
$ws.Range("E28").Value = "  +3.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "23.01"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "162.81"
$ws.Range("E30").Value = "  -5.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.67"
$ws.Range("E31").Value = "  -5.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0875"
$ws.Range("E32").Value = "  -1.39%  "
$ws.Range("E33").Value = "  +1.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.82"
$ws.Range("E34").Value = "  -3.70%  "
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("E36").Value = "  -7.71%  "
$ws.Range("E37").Value = "  -1.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.88"
$ws.Range("E38").Value = "  +8.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0347"
$ws.Range("E39").Value = "  -3.77%  "
$ws.Range("E40").Value = "  -7.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "103.28"
$ws.Range("E41").Value = "  +21.29%  "
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "69.56"
$ws.Range("E43").Value = "  -1.45%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.225"
$ws.Range("E45").Value = "  -4.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "114.89"
$ws.Range("E46").Value = "  +2.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.95"
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.99"
$ws.Range("E48").Value = "  -1.62%  "
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.80"
$ws.Range("E49").Value = "  +3.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.28"
$ws.Range("E50").Value = "  -2.93%  "
$ws.Range("E51").Value = "  -1.14%  "
